# Active Learning on the Cloud - PT4.pptx
# - Move the "Picture 65" image (slide 3) up slightly.
# - Re-route the elbow connector ("Connector: Elbow 90", slide 3) that starts
#   on that picture, so its geometry matches the picture's new position.
# - Merge the two runs "Performance " + "requirement" on slide 8 into a
#   single run (keeping the second run's formatting/dirty flag).

function Find-ShapeByName($shapes, $name) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -eq $name) { return $sh }
    }
    return $null
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 3: reposition the picture and its attached connector.
# ---------------------------------------------------------------------
$slide3 = $p.Slides.Item(3)

# "Picture 65" - only the vertical offset changes (x stays 5269909 EMU).
$picture = Find-ShapeByName $slide3.Shapes "Picture 65"
if ($picture -eq $null) { $picture = $slide3.Shapes.Item(15) }
$picture.Top = 437.95191956377954   # -> 5561989 EMU

# "Connector: Elbow 90" - connected (stCxn) to the picture above; its frame
# is stored explicitly in the XML and is not recomputed automatically, so
# update it to match the picture's new location.
$connector = Find-ShapeByName $slide3.Shapes "Connector: Elbow 90"
if ($connector -eq $null) { $connector = $slide3.Shapes.Item(23) }
$connector.Left = 428.8300934401575    # -> 5446142 EMU
$connector.Top = 428.9958648716535     # -> 5448247 EMU
$connector.Width = 17.91224392440945   # -> 227485 EMU

# ---------------------------------------------------------------------
# Slide 8: merge "Performance " + "requirement" into a single text run.
# ---------------------------------------------------------------------
$slide8 = $p.Slides.Item(8)
$contentShape = Find-ShapeByName $slide8.Shapes "Content Placeholder 2"
if ($contentShape -eq $null) { $contentShape = $slide8.Shapes.Item(2) }
$tr = $contentShape.TextFrame.TextRange

$splitText = "Performance "
$paraCount = $tr.Count
for ($i = 1; $i -le $paraCount; $i++) {
    $para = $tr.Paragraphs($i, 1)
    $paraText = $para.Text.TrimEnd("`r")
    if ($paraText -eq "Performance requirement") {
        $firstRun = $tr.Characters($para.Start, $splitText.Length)
        $firstRun.Text = ""
        $secondRun = $tr.Characters($para.Start, "requirement".Length)
        $secondRun.InsertBefore($splitText) | Out-Null
        break
    }
}
